$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; existing rows (header + data) shift down by one.
$ws.Rows("1:1").Insert()

# New title row (row 1) - plain text, no header styling.
$ws.Range("A1").Value = "2023年采购数据"

# Two new data rows appended at the bottom (rows 6 and 7).
$ws.Range("A6").Value = "D101"
$ws.Range("B6").Value = 400
$ws.Range("C6").Value = "Supplier D"
$ws.Range("D6").Value = 4000
$ws.Range("E6").Value = "M"

$ws.Range("A7").Value = "E202"
$ws.Range("B7").Value = 500
$ws.Range("C7").Value = "Supplier E"
$ws.Range("D7").Value = 5000
$ws.Range("E7").Value = "N"
